# Generate Report for Handback
# Fill in the Correspond Handoff/Handback datetimes for the second
# (d3f61a6d...) row of each language sheet, which previously shared the
# same placeholder values as the first row.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-11 09:36:10"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "2016-03-11 09:36:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-11 09:36:13"
$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "2016-03-11 09:36:34"
